# events_selection.xlsx: add the "Book of the Month Giveaway" rows to the
# Events list (one instance "This" Monday, one instance "Next" Tuesday),
# nudge the saved selection, and pull in the small (8pt) font that Excel
# registers alongside the sheet's phonetic-guide info for the new text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- introduce the 8pt Calibri font (mirrors the one Excel adds for the
#     phonetic-guide font table) without leaving any visible cell touched:
#     flip an existing, already-centred data cell to size 8 and straight
#     back to its original size. ---
$fontProbe = $ws.Range("B2")
$originalFontSize = $fontProbe.Font.Size()
$fontProbe.Font.Size = 8
$fontProbe.Font.Size = $originalFontSize

# --- row 18: Book of the Month Giveaway, This / Monday, 12:00 AM ---
$ws.Range("A18").Value = "Yes"
$ws.Range("A18").HorizontalAlignment = -4108
$ws.Range("B18").Value = "Book of the Month Giveaway: The Art of Keeping Secrets"
$ws.Range("C18").Value = "Book"
$ws.Range("D18").Value = "This"
$ws.Range("E18").Value = "Monday"
$ws.Range("F18").Value = 0
$ws.Range("F18").NumberFormat = "[$-409]h:mm\ AM/PM;@"

# --- row 19: Book of the Month Giveaway, Next / Tuesday, 1:00 AM ---
$ws.Range("A19").Value = "Yes"
$ws.Range("A19").HorizontalAlignment = -4108
$ws.Range("B19").Value = "Book of the Month Giveaway: The Art of Keeping Secrets"
$ws.Range("C19").Value = "Book"
$ws.Range("D19").Value = "Next"
$ws.Range("E19").Value = "Tuesday"
$ws.Range("F19").Value = 0.041666666666666699
$ws.Range("F19").NumberFormat = "[$-409]h:mm\ AM/PM;@"

# --- match the saved selection left behind in the workbook ---
$ws.Range("D20").Select()

Write-Output "Added Book of the Month Giveaway rows (18-19)."
